$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -8.757999999999996
$ws.Range("A4").Value = -22.34640000000001
$ws.Range("B4").Value = 5.4985
$ws.Range("A6").Value = -22.75630000000001
$ws.Range("A7").Value = -19.82699999999998
$ws.Range("C7").Value = -12.95160000000001
$ws.Range("C8").Value = -12.36279999999999
$ws.Range("B9").Value = 6.350300000000001
$ws.Range("C10").Value = -13.71929999999999
$ws.Range("B12").Value = 4.800099999999997
$ws.Range("C13").Value = -13.4166
$ws.Range("D13").Value = -8.787699999999996
$ws.Range("A16").Value = -21.87520000000001
$ws.Range("C16").Value = -12.5765
$ws.Range("B17").Value = 5.468599999999999
$ws.Range("B18").Value = 6.441099999999993
$ws.Range("B19").Value = 9.269199999999998
$ws.Range("A20").Value = -22.10660000000001
$ws.Range("B20").Value = 5.721999999999998
$ws.Range("D20").Value = -7.941400000000008
$ws.Range("D25").Value = -7.586399999999994
$ws.Range("B26").Value = 4.358300000000005
$ws.Range("A28").Value = -21.90619999999999
$ws.Range("A29").Value = -21.24109999999997
$ws.Range("C30").Value = -11.9635
$ws.Range("B31").Value = 4.057399999999999
$ws.Range("A32").Value = -21.29420000000002
$ws.Range("D34").Value = -7.493300000000001
$ws.Range("B39").Value = 9.249100000000004
$ws.Range("D39").Value = -8.102899999999993
$ws.Range("A40").Value = -21.77019999999998
$ws.Range("B40").Value = 6.031599999999998
$ws.Range("C40").Value = -12.5336
$ws.Range("B41").Value = 9.708199999999984
$ws.Range("B42").Value = 9.59119999999999
$ws.Range("B43").Value = 6.329400000000003
$ws.Range("C44").Value = -13.6155
$ws.Range("A46").Value = -21.79
$ws.Range("B47").Value = 5.282999999999999
$ws.Range("B48").Value = 5.501300000000003
$ws.Range("A51").Value = -22.0711
$ws.Range("D51").Value = -8.594900000000001
$ws.Range("A52").Value = -22.13049999999999
$ws.Range("A57").Value = -21.85600000000001
$ws.Range("A59").Value = -22.1855
$ws.Range("D59").Value = -8.330499999999995
$ws.Range("D61").Value = -8.225799999999994
$ws.Range("A62").Value = -22.13600000000002
$ws.Range("B63").Value = 4.816699999999998
$ws.Range("B64").Value = 5.394200000000003
$ws.Range("D64").Value = -6.981799999999996
$ws.Range("A66").Value = -21.39690000000001
$ws.Range("A73").Value = -20.2603
$ws.Range("A74").Value = -21.68549999999998
$ws.Range("B76").Value = 5.7893
$ws.Range("D78").Value = -8.247000000000003
$ws.Range("B81").Value = 5.117400000000005
$ws.Range("D83").Value = -8.608700000000004
$ws.Range("B89").Value = 5.3767
$ws.Range("C89").Value = -14.41570000000001
$ws.Range("C91").Value = -12.60360000000001
$ws.Range("A92").Value = -21.60950000000002
$ws.Range("D92").Value = -6.134600000000002
$ws.Range("B94").Value = 4.762999999999993
$ws.Range("D98").Value = -7.1249
$ws.Range("A100").Value = -22.121
$ws.Range("D100").Value = -8.204199999999993
